$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells retain their original text formatting,
# since many price values look numeric (e.g. "1.003", "31.316.73") and
# would otherwise be auto-converted/rounded by Excel when assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '31.316.73'
$ws.Range("E2").Value = '  +3.02%  '
$ws.Range("D3").Value = '2.007.00'
$ws.Range("E3").Value = '  +7.20%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = '0.7801'
$ws.Range("E5").Value = '  +65.79%  '
$ws.Range("D6").Value = '259.73'
$ws.Range("E6").Value = '  +6.20%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = '0.3582'
$ws.Range("E8").Value = '  +24.58%  '
$ws.Range("D9").Value = '28.43'
$ws.Range("E9").Value = '  +29.98%  '
$ws.Range("D10").Value = '0.07067'
$ws.Range("E10").Value = '  +8.88%  '
$ws.Range("D11").Value = '0.8593'
$ws.Range("E11").Value = '  +17.58%  '
$ws.Range("D12").Value = '0.08200'
$ws.Range("E12").Value = '  +5.29%  '
$ws.Range("D13").Value = '2.007.09'
$ws.Range("E13").Value = '  +7.20%  '
$ws.Range("D14").Value = '101.44'
$ws.Range("E14").Value = '  +1.70%  '
$ws.Range("D15").Value = '5.628'
$ws.Range("E15").Value = '  +8.83%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '15.40'
$ws.Range("E16").Value = '  +17.46%  '
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").Value = '274.62'
$ws.Range("E17").Value = '  -4.06%  '
$ws.Range("D18").Value = '31.311.25'
$ws.Range("E18").Value = '  +3.05%  '
$ws.Range("D19").Value = '5.949'
$ws.Range("E19").Value = '  +11.94%  '
$ws.Range("D20").Value = '0.000008002'
$ws.Range("E20").Value = '  +6.81%  '
$ws.Range("D21").Value = '2.272.26'
$ws.Range("E21").Value = '  +7.60%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").Value = '7.180'
$ws.Range("E24").Value = '  +13.58%  '
$ws.Range("D25").Value = '10.09'
$ws.Range("E25").Value = '  +11.48%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '165.92'
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = '0.1496'
$ws.Range("E27").Value = '  +54.75%  '
$ws.Range("D28").Value = '20.06'
$ws.Range("E28").Value = '  +6.01%  '
$ws.Range("D29").Value = '2.393'
$ws.Range("E29").Value = '  +26.24%  '
$ws.Range("E30").Value = '  +9.32%  '
$ws.Range("D31").Value = '4.642'
$ws.Range("E31").Value = '  +9.70%  '
$ws.Range("D32").Value = '1.360'
$ws.Range("E32").Value = '  +3.20%  '
$ws.Range("D33").Value = '4.435'
$ws.Range("E33").Value = '  +6.97%  '
$ws.Range("D34").Value = '0.05237'
$ws.Range("E34").Value = '  +8.98%  '
$ws.Range("D35").Value = '0.7794'
$ws.Range("E35").Value = '  +13.15%  '
$ws.Range("D36").Value = '1.222'
$ws.Range("E36").Value = '  +8.50%  '
$ws.Range("D37").Value = '2.815'
$ws.Range("E37").Value = '  +3.24%  '
$ws.Range("D38").Value = '0.02007'
$ws.Range("E38").Value = '  +5.65%  '
$ws.Range("D39").Value = '2.945'
$ws.Range("E39").Value = '  +3.60%  '
$ws.Range("D40").Value = '6.718'
$ws.Range("E40").Value = '  +7.02%  '
$ws.Range("D41").Value = '80.24'
$ws.Range("E41").Value = '  +5.71%  '
$ws.Range("D42").Value = '0.4748'
$ws.Range("E42").Value = '  +12.30%  '
$ws.Range("D43").Value = '2.157'
$ws.Range("E43").Value = '  +9.80%  '
$ws.Range("D44").Value = '107.58'
$ws.Range("E44").Value = '  +6.34%  '
$ws.Range("D45").Value = '0.8592'
$ws.Range("E45").Value = '  +4.27%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.821'
$ws.Range("E47").Value = '  +11.49%  '
$ws.Range("D48").Value = '9.998'
$ws.Range("E48").Value = '  +2.45%  '
$ws.Range("D49").Value = '0.4364'
$ws.Range("E49").Value = '  +11.62%  '
$ws.Range("D50").Value = '36.94'
$ws.Range("E50").Value = '  +5.52%  '
$ws.Range("D51").Value = '0.1201'
$ws.Range("E51").Value = '  +14.45%  '
